$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting the existing quarterly
# columns (D:K) right to (E:L). This matches the author's addition of a
# new latest-quarter column while keeping the 8 prior quarters.
$ws.Columns("D:D").Insert()

# The freshly inserted column D has no number formatting yet (it picks up
# the plain default style). Copy the formatting from the column that used
# to be D (now shifted to E) back onto D so every row gets the right
# style (date style on the "Period Ending" rows, #,##0 style elsewhere).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest-quarter figures.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 1500
$ws.Range("D9").Value = 0
$ws.Range("D10").Value = 1500
$ws.Range("D12").Value = 22000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 29900
$ws.Range("D18").Value = -28400
$ws.Range("D20").Value = 1200
$ws.Range("D21").Value = -27100
$ws.Range("D22").Value = 2300
$ws.Range("D23").Value = -29600
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -29600
$ws.Range("D27").Value = -29600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1200
$ws.Range("D33").Value = -29600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -29600
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 52600
$ws.Range("D42").Value = 67700
$ws.Range("D43").Value = 3400
$ws.Range("D44").Value = 800
$ws.Range("D45").Value = 4000
$ws.Range("D46").Value = 128500
$ws.Range("D47").Value = 29200
$ws.Range("D48").Value = 9200
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 168300
$ws.Range("D57").Value = 7900
$ws.Range("D58").Value = 30700
$ws.Range("D59").Value = 25700
$ws.Range("D60").Value = 64400
$ws.Range("D61").Value = 31000
$ws.Range("D62").Value = 100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 95400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -704500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 72900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -29600
$ws.Range("D83").Value = 200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -35000
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -9500
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 63900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 19400
